$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows right after the header row (before current row 2),
# shifting existing data rows down.
$insertRange = $ws.Range("A2:C7")
$insertRange.EntireRow.Insert()
$insertRange.EntireRow.ClearFormats()

$newTopRows = @(
    @(-3.130717563629151, 5.486354422569275, -2.054217553138733),
    @(-3.261763083934784, 5.438373637199402, -2.230073320865632),
    @(-3.289464282989502, 5.444673538208008, -2.207874870300293),
    @(-3.29152911901474,  5.474263513088226, -1.988579791784286),
    @(-3.39021909236908,  5.475549221038817, -1.85912013053894),
    @(-3.51887332201004,  5.510936594009399, -1.821529471874237)
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Append 4 new rows at the end (rows 28-31 after the insert above).
$newBottomRows = @(
    @(1.438675880432129, 5.703988456726075, 1.146768474578857),
    @(1.517132639884949, 5.866671967506409, 1.15747617483139),
    @(1.545208883285523, 5.865323352813721, 1.205629134178161),
    @(1.573752522468567, 5.771291553974152, 1.240318953990936)
)

$r = 28
foreach ($row in $newBottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
